$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $rowIndex, $colIndex, $value) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item($colIndex)
    $r = $cell.Range
    $r.Text = $value
    $r.Font.Bold = 1
    $r.Font.Size = 12
    $r.Font.SizeBi = 12
}

# RETENTION section -> "Ratio" row -> value cell
Set-CellValue $t 24 2 "0.6"

# QUESTION AND ANSWER TASK -> Answer Recall Lenient (ARL) -> value cell
Set-CellValue $t 44 2 "0.25"

# QUESTION AND ANSWER TASK -> Answer Recall Strict (ARS) -> value cell
Set-CellValue $t 45 2 "0.1666"

# QUESTION AND ANSWER TASK -> Answer Recall Average (ARA) -> value cell
Set-CellValue $t 46 2 "0.2083"
